$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 2852
$ws.Range("D4").Value = 632
$ws.Range("D5").Value = 38
$ws.Range("D6").Value = 16

$ws.Range("D2:D6").Interior.Pattern = -4142

$ws.Range("A12").Select()
